$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11: fill in the 1/13 diary entry (date, date+time, and new text values)
# (Match the number formats used by the row above it, row 10, so the
#  values are stored as real numbers rather than text.)
$ws.Range("A11").NumberFormat = "m/d"
$ws.Range("A11").Value = 43843
$ws.Range("B11").NumberFormat = "h:mmAM/PM"
$ws.Range("B11").Value = 43843.91875
$ws.Range("C11").Value = "n/a"
$ws.Range("D11").Value = "Figure out how to keep local repo and local directory up-to-date"
$ws.Range("E11").Value = "Solution: navigate to local directory via terminal, perform — git pull upstream master"
$ws.Range("F11").Value = "I will definitely run into other issues on git and GitHub as I go through this course and find solutions on the go"
$ws.Range("G11").Value = "Hopeful; Hoping that I become more capable by the end of this course than I was in the beginning "

# Row 12: clear placeholder text, leaving the row blank
$ws.Range("A12:G12").ClearContents()

# Row 13: clear placeholder text in column A only (B13:G13 were already blank)
$ws.Range("A13").ClearContents()
